$wb = $excel.ActiveWorkbook

# --- interfaceOperations sheet (sheet2): full method list incl. inherited Object methods ---
$ws2 = $wb.Worksheets.Item("interfaceOperations")
$ws2.Cells.Item(1,1).Value = "Interface Name"
$ws2.Cells.Item(1,2).Value = "Operation Signature"
$ws2.Cells.Item(1,3).Value = "Operation Modifier"
$ws2.Cells.Item(1,4).Value = "Return Type"
$ws2.Cells.Item(2,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(2,2).Value = "equals(java.lang.Object)"
$ws2.Cells.Item(2,3).Value = "public"
$ws2.Cells.Item(2,4).Value = "boolean"
$ws2.Cells.Item(3,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(3,2).Value = "GatewayController()"
$ws2.Cells.Item(3,3).Value = "public"
$ws2.Cells.Item(3,4).Value = "void"
$ws2.Cells.Item(4,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(4,2).Value = "toString()"
$ws2.Cells.Item(4,3).Value = "public"
$ws2.Cells.Item(4,4).Value = "java.lang.String"
$ws2.Cells.Item(5,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(5,2).Value = "getClass()"
$ws2.Cells.Item(5,3).Value = "public"
$ws2.Cells.Item(5,4).Value = "java.lang.Class"
$ws2.Cells.Item(6,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(6,2).Value = "notifyAll()"
$ws2.Cells.Item(6,3).Value = "public"
$ws2.Cells.Item(6,4).Value = "void"
$ws2.Cells.Item(7,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(7,2).Value = "hashCode()"
$ws2.Cells.Item(7,3).Value = "public"
$ws2.Cells.Item(7,4).Value = "int"
$ws2.Cells.Item(8,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(8,2).Value = "wait()"
$ws2.Cells.Item(8,3).Value = "public"
$ws2.Cells.Item(8,4).Value = "void"
$ws2.Cells.Item(9,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(9,2).Value = "hi()"
$ws2.Cells.Item(9,3).Value = "public"
$ws2.Cells.Item(9,4).Value = "reactor.core.publisher.Mono"
$ws2.Cells.Item(10,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(10,2).Value = "notify()"
$ws2.Cells.Item(10,3).Value = "public"
$ws2.Cells.Item(10,4).Value = "void"
$ws2.Cells.Item(11,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(11,2).Value = "wait(long)"
$ws2.Cells.Item(11,3).Value = "public"
$ws2.Cells.Item(11,4).Value = "void"
$ws2.Cells.Item(12,1).Value = "org.andante.GatewayController"
$ws2.Cells.Item(12,2).Value = "wait(long, int)"
$ws2.Cells.Item(12,3).Value = "public"
$ws2.Cells.Item(12,4).Value = "void"

# --- methodNumberOfLines sheet (sheet11): add constructors/lambdas, reorder ---
$ws11 = $wb.Worksheets.Item("methodNumberOfLines")

# Column C holds numeric-looking text ("1", "3", ...). Force text storage (matches
# the source data, which stores line counts as shared strings, not numbers) by
# switching the column to Text format before writing, then clearing the format back
# to the default style afterwards (keeps the text type, drops the custom number format).
$ws11.Columns.Item(3).NumberFormat = "@"

$ws11.Cells.Item(1,1).Value = "Class Name"
$ws11.Cells.Item(1,2).Value = "Method Signature"
$ws11.Cells.Item(1,3).Value = "Number of Lines"
$ws11.Cells.Item(2,1).Value = "org.andante.config.security.role.KeycloakRole"
$ws11.Cells.Item(2,2).Value = "values()"
$ws11.Cells.Item(2,3).Value = "1"
$ws11.Cells.Item(3,1).Value = "org.andante.config.security.role.KeycloakRole"
$ws11.Cells.Item(3,2).Value = "valueOf(java.lang.String)"
$ws11.Cells.Item(3,3).Value = "1"
$ws11.Cells.Item(4,1).Value = "org.andante.config.security.role.KeycloakRole"
$ws11.Cells.Item(4,2).Value = "KeycloakRole(java.lang.String, int, java.lang.String)"
$ws11.Cells.Item(4,3).Value = "3"
$ws11.Cells.Item(5,1).Value = "org.andante.config.security.role.KeycloakRole"
$ws11.Cells.Item(5,2).Value = "getName()"
$ws11.Cells.Item(5,3).Value = "3"
$ws11.Cells.Item(6,1).Value = "org.andante.config.security.role.KeycloakRole"
$ws11.Cells.Item(6,2).Value = "`$values()"
$ws11.Cells.Item(6,3).Value = "1"
$ws11.Cells.Item(7,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(7,2).Value = "gatewayRoutes(org.springframework.cloud.gateway.route.builder.RouteLocatorBuilder)"
$ws11.Cells.Item(7,3).Value = "3"
$ws11.Cells.Item(8,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(8,2).Value = "GatewayConfiguration(org.springframework.cloud.gateway.filter.factory.TokenRelayGatewayFilterFactory)"
$ws11.Cells.Item(8,3).Value = "3"
$ws11.Cells.Item(9,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(9,2).Value = "lambda`$gatewayRoutes`$11(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(9,3).Value = "1"
$ws11.Cells.Item(10,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(10,2).Value = "lambda`$gatewayRoutes`$10(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(10,3).Value = "1"
$ws11.Cells.Item(11,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(11,2).Value = "lambda`$gatewayRoutes`$9(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(11,3).Value = "1"
$ws11.Cells.Item(12,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(12,2).Value = "lambda`$gatewayRoutes`$8(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(12,3).Value = "1"
$ws11.Cells.Item(13,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(13,2).Value = "lambda`$gatewayRoutes`$7(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(13,3).Value = "1"
$ws11.Cells.Item(14,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(14,2).Value = "lambda`$gatewayRoutes`$6(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(14,3).Value = "1"
$ws11.Cells.Item(15,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(15,2).Value = "lambda`$gatewayRoutes`$5(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(15,3).Value = "1"
$ws11.Cells.Item(16,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(16,2).Value = "lambda`$gatewayRoutes`$4(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(16,3).Value = "1"
$ws11.Cells.Item(17,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(17,2).Value = "lambda`$gatewayRoutes`$3(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(17,3).Value = "1"
$ws11.Cells.Item(18,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(18,2).Value = "lambda`$gatewayRoutes`$2(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(18,3).Value = "1"
$ws11.Cells.Item(19,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(19,2).Value = "lambda`$gatewayRoutes`$1(org.springframework.cloud.gateway.route.builder.PredicateSpec)"
$ws11.Cells.Item(19,3).Value = "1"
$ws11.Cells.Item(20,1).Value = "org.andante.config.gateway.GatewayConfiguration"
$ws11.Cells.Item(20,2).Value = "lambda`$gatewayRoutes`$0(org.springframework.cloud.gateway.route.builder.GatewayFilterSpec)"
$ws11.Cells.Item(20,3).Value = "1"
$ws11.Cells.Item(21,1).Value = "org.andante.config.security.filter.CrossOriginRequestSharingFilter"
$ws11.Cells.Item(21,2).Value = "CrossOriginRequestSharingFilter()"
$ws11.Cells.Item(21,3).Value = "6"
$ws11.Cells.Item(22,1).Value = "org.andante.config.security.filter.CrossOriginRequestSharingFilter"
$ws11.Cells.Item(22,2).Value = "filter(org.springframework.web.server.ServerWebExchange, org.springframework.web.server.WebFilterChain)"
$ws11.Cells.Item(22,3).Value = "11"
$ws11.Cells.Item(23,1).Value = "org.andante.GatewayApplication"
$ws11.Cells.Item(23,2).Value = "GatewayApplication()"
$ws11.Cells.Item(23,3).Value = "1"
$ws11.Cells.Item(24,1).Value = "org.andante.GatewayApplication"
$ws11.Cells.Item(24,2).Value = "main(java.lang.String[])"
$ws11.Cells.Item(24,3).Value = "3"
$ws11.Cells.Item(25,1).Value = "org.andante.GatewayController"
$ws11.Cells.Item(25,2).Value = "GatewayController()"
$ws11.Cells.Item(25,3).Value = "1"
$ws11.Cells.Item(26,1).Value = "org.andante.GatewayController"
$ws11.Cells.Item(26,2).Value = "hi()"
$ws11.Cells.Item(26,3).Value = "3"
$ws11.Cells.Item(27,1).Value = "org.andante.config.security.converter.KeycloakRealmRoleConverter"
$ws11.Cells.Item(27,2).Value = "KeycloakRealmRoleConverter()"
$ws11.Cells.Item(27,3).Value = "1"
$ws11.Cells.Item(28,1).Value = "org.andante.config.security.converter.KeycloakRealmRoleConverter"
$ws11.Cells.Item(28,2).Value = "convert(org.springframework.security.oauth2.jwt.Jwt)"
$ws11.Cells.Item(28,3).Value = "8"
$ws11.Cells.Item(29,1).Value = "org.andante.config.security.converter.KeycloakRealmRoleConverter"
$ws11.Cells.Item(29,2).Value = "convert(java.lang.Object)"
$ws11.Cells.Item(29,3).Value = "1"
$ws11.Cells.Item(30,1).Value = "org.andante.config.security.converter.KeycloakRealmRoleConverter"
$ws11.Cells.Item(30,2).Value = "lambda`$convert`$0(java.lang.String)"
$ws11.Cells.Item(30,3).Value = "1"
$ws11.Cells.Item(31,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(31,2).Value = "securityWebFilterChain(org.springframework.security.config.web.server.ServerHttpSecurity)"
$ws11.Cells.Item(31,3).Value = "4"
$ws11.Cells.Item(32,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(32,2).Value = "corsConfiguration()"
$ws11.Cells.Item(32,3).Value = "14"
$ws11.Cells.Item(33,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(33,2).Value = "jwtDecoder()"
$ws11.Cells.Item(33,3).Value = "3"
$ws11.Cells.Item(34,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(34,2).Value = "jwtAuthenticationConverter()"
$ws11.Cells.Item(34,3).Value = "5"
$ws11.Cells.Item(35,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(35,2).Value = "SecurityConfiguration(org.andante.config.security.converter.KeycloakRealmRoleConverter)"
$ws11.Cells.Item(35,3).Value = "8"
$ws11.Cells.Item(36,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(36,2).Value = "lambda`$securityWebFilterChain`$2(org.springframework.security.config.web.server.ServerHttpSecurity`$OAuth2ResourceServerSpec)"
$ws11.Cells.Item(36,3).Value = "1"
$ws11.Cells.Item(37,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(37,2).Value = "lambda`$securityWebFilterChain`$1(org.springframework.security.config.web.server.ServerHttpSecurity`$OAuth2ResourceServerSpec`$JwtSpec)"
$ws11.Cells.Item(37,3).Value = "1"
$ws11.Cells.Item(38,1).Value = "org.andante.config.security.SecurityConfiguration"
$ws11.Cells.Item(38,2).Value = "lambda`$securityWebFilterChain`$0(int)"
$ws11.Cells.Item(38,3).Value = "1"

$ws11.Columns.Item(3).ClearFormats()
